$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("A1").Value = "ASIN"
$ws.Range("B1").Value = "changepoint_prior_scale"
$ws.Range("C1").Value = "seasonality_prior_scale"
$ws.Range("D1").Value = "holidays_prior_scale"
$ws.Range("E1").Value = "RMSE_Mean"
$ws.Range("F1").Value = "RMSE_P70"
$ws.Range("G1").Value = "RMSE_P80"
$ws.Range("H1").Value = "RMSE_P90"
$ws.Range("I1").Value = "Total Tests"

# Data row values
$ws.Range("A2").Value = "B08F7BHDLY"
$ws.Range("B2").Value = 0.5
$ws.Range("C2").Value = 0.05
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 22.65364209128413
$ws.Range("F2").Value = 12.05456345124119
$ws.Range("G2").Value = 14.72667647502314
$ws.Range("H2").Value = 33.89874628950162
$ws.Range("I2").Value = 126

# Build the header style on a single cell first (so only one new style entry
# gets created), then copy that formatting across the rest of the header row.
$headerCell = $ws.Range("A1")
$headerCell.Borders.LineStyle = 1
$headerCell.Borders.Weight = 2
$headerCell.Font.Bold = $true
$headerCell.HorizontalAlignment = -4108
$headerCell.VerticalAlignment = -4160

$headerCell.Copy()
$ws.Range("A1:I1").PasteSpecial(-4122)
